$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current last data row (118), pushing the
# existing row 118 down to row 120 (its contents/formatting move with it).
$ws.Rows("118:119").Insert()

# New row 118: Packham's Triumph, "Primera" quality, week of 2021-09-09.
$ws.Range("A118").Value = 4
$ws.Range("B118").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C118").Value = "Los Lagos"
$ws.Range("D118").Value = 44448
$ws.Range("E118").Value = 10
$ws.Range("F118").Value = "Fruta"
$ws.Range("G118").Value = 100104
$ws.Range("H118").Value = "Frutos de pepita"
$ws.Range("I118").Value = 100104005
$ws.Range("J118").Value = "Pera"
$ws.Range("K118").Value = "Packham's Triumph"
$ws.Range("L118").Value = "Primera"
$ws.Range("M118").Value = 200
$ws.Range("N118").Value = 16000
$ws.Range("O118").Value = 16000
$ws.Range("P118").Value = 16000
$ws.Range("Q118").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R118").Value = "Región de O'Higgins"
$ws.Range("S118").Value = 1067
$ws.Range("T118").Value = 15

# New row 119: Packham's Triumph, "Segunda" quality, same week.
$ws.Range("A119").Value = 4
$ws.Range("B119").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C119").Value = "Los Lagos"
$ws.Range("D119").Value = 44448
$ws.Range("E119").Value = 10
$ws.Range("F119").Value = "Fruta"
$ws.Range("G119").Value = 100104
$ws.Range("H119").Value = "Frutos de pepita"
$ws.Range("I119").Value = 100104005
$ws.Range("J119").Value = "Pera"
$ws.Range("K119").Value = "Packham's Triumph"
$ws.Range("L119").Value = "Segunda"
$ws.Range("M119").Value = 100
$ws.Range("N119").Value = 11000
$ws.Range("O119").Value = 11000
$ws.Range("P119").Value = 11000
$ws.Range("Q119").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R119").Value = "Región de O'Higgins"
$ws.Range("S119").Value = 733
$ws.Range("T119").Value = 15
